$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, pushing the existing rows 119-177 down to 120-178.
$ws.Rows("119:119").Insert()

# Populate the newly inserted row 119 with the new data point.
$ws.Cells.Item(119, 1).Value = 11
$ws.Cells.Item(119, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(119, 3).Value = "Bíobío"
$ws.Cells.Item(119, 4).Value = 44609
$ws.Cells.Item(119, 5).Value = 8
$ws.Cells.Item(119, 6).Value = 100114001
$ws.Cells.Item(119, 7).Value = "Papa"
$ws.Cells.Item(119, 8).Value = "Asterix"
$ws.Cells.Item(119, 9).Value = "1a (cosecha)"
$ws.Cells.Item(119, 10).Value = 4000
$ws.Cells.Item(119, 11).Value = 7500
$ws.Cells.Item(119, 12).Value = 8000
$ws.Cells.Item(119, 13).Value = 7750
$ws.Cells.Item(119, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(119, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(119, 16).Value = 310
$ws.Cells.Item(119, 17).Value = 25
$ws.Cells.Item(119, 18).Value = "Hortaliza"
